# The data rows for this subset (rows 2, 4, 5, 6 on Sheet1) need to be
# re-ordered: row 2 <-> row 6 and row 4 <-> row 5 swap their contents
# (row 3 is left untouched). We swap the columns that actually differ
# between the rows: D (Fecha), I (Calidad), J (Volumen), K (Precio
# minimo), L (Precio maximo), M (Precio promedio ponderado), N (Unidad
# de comercializacion), P (Precio $/Kg) and Q (Kg o Unidades).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("D", "I", "J", "K", "L", "M", "N", "P", "Q")

function Swap-RowValues($ws, $rowA, $rowB, $cols) {
    foreach ($col in $cols) {
        $cellA = $ws.Range($col + $rowA)
        $cellB = $ws.Range($col + $rowB)
        $valA = $cellA.Value2
        $valB = $cellB.Value2
        $cellA.Value = $valB
        $cellB.Value = $valA
    }
}

Swap-RowValues $ws 2 6 $cols
Swap-RowValues $ws 4 5 $cols
